# Insert a new "Disciplinary Procedures" agenda item before the existing
# "Renovation in Summer" item, and move the "_GoBack" bookmark (which sits
# right after "Admission" in the original) into the new item's text, split
# between "Disciplinar" and "y Procedures" - matching how Word leaves the
# _GoBack bookmark at the most recent edit point.

$d = $word.ActiveDocument

# Locate the "Renovation in Summer" bullet via Find.
$found = $d.Content
$found.Find.Execute("Renovation in Summer", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fStart = $found.Start
$fEnd = $found.End

# Resolve the paragraph object that contains the found text.
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $fStart -and $pp.Range.End -ge $fEnd) {
        $idx = $i
        break
    }
}

$renovPara = $d.Paragraphs.Item($idx)

# Insert a new empty list paragraph right before it; it inherits the same
# pPr (ListParagraph style / numPr ilvl=1 numId=1).
$renovPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($idx)
$newPara.Range.Text = "Disciplinary Procedures"

# Move the "_GoBack" bookmark out of the old "Admission" paragraph and into
# the new paragraph, right after "Disciplinar" (11 characters in), so the
# run is split into "Disciplinar" / bookmark / "y Procedures".
$bmPos = $newPara.Range.Start + 11
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
